$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999434018183
$ws.Range("A2").Value = 0.99872453683934748
$ws.Range("A3").Value = 0.99710000932592568
$ws.Range("A4").Value = 1.0017670985283256
$ws.Range("A5").Value = 0.99334510500897233
$ws.Range("A6").Value = 0.97541515191171579
$ws.Range("A7").Value = 0.97323035028612925
$ws.Range("A8").Value = 0.96976597915979446
$ws.Range("A9").Value = 0.96823358908851087
$ws.Range("A10").Value = 0.96763945003106833
$ws.Range("A11").Value = 0.96758307054981629
$ws.Range("A12").Value = 0.96773241808385324
$ws.Range("A13").Value = 0.97132087070902251
$ws.Range("A14").Value = 0.97408886198162925
$ws.Range("A15").Value = 0.9777792683766563
$ws.Range("A16").Value = 0.9752732430363964
$ws.Range("A17").Value = 0.9715657680343458
$ws.Range("A18").Value = 0.97045687345367904
$ws.Range("A19").Value = 0.99764436283918378
$ws.Range("A20").Value = 0.99052736647185635
$ws.Range("A21").Value = 0.98912887916073822
$ws.Range("A22").Value = 0.98786437605483224
$ws.Range("A23").Value = 0.99006607522045764
$ws.Range("A24").Value = 0.97704608080685751
$ws.Range("A25").Value = 0.97058923093234495
$ws.Range("A26").Value = 0.97298748410395186
$ws.Range("A27").Value = 0.96985208242958976
$ws.Range("A28").Value = 0.95908294654280257
$ws.Range("A29").Value = 0.95174753757021158
$ws.Range("A30").Value = 0.94811541438038316
$ws.Range("A31").Value = 0.95127671694848792
$ws.Range("A32").Value = 0.94959744816381431
$ws.Range("A33").Value = 0.94907743691377489
